$d = $word.ActiveDocument

# Replace each manual line break with a single space, joining the
# separate <w:t> runs that were split by <w:br/> elements into one
# continuous run of text.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2)

# The very last paragraph mark was preceded by a trailing manual line
# break (with no text after it); the bulk replace above turned that
# into a trailing space that should not exist in the final text.
$endRange = $d.Content
$endRange.Collapse(0)
[void]$endRange.MoveEnd(1, -1)
[void]$endRange.MoveStart(1, -1)
if ($endRange.Text -eq " ") {
    $endRange.Text = ""
}
